# Updates the crypto price/volume table (and two rank-order swaps) to
# match the latest scrape, per commit "Updated cryptos list on Wed Jun 21
# 15:40:34 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column (D) cells whose new value would otherwise be auto-parsed by
# Excel as a number (stripping the trailing zero / changing the stored
# type). Forcing Text format first keeps them stored as literal strings,
# exactly like the rest of the (locale-formatted) price column.
$textPriceCells = @(
    'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D11', 'D12', 'D13', 'D14', 'D15', 
    'D17', 'D18', 'D19', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 
    'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 
    'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D48', 
    'D49', 'D50', 'D51'
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# New cell values: coin name / link / price / 1h-volume-change columns.
$updates = @{
    'D2' = '29.805.67'
    'E2' = '  +11.18%  '
    'D3' = '1.845.17'
    'E3' = '  +7.19%  '
    'D4' = '0.9974'
    'E4' = '  -0.18%  '
    'D5' = '248.77'
    'E5' = '  +3.76%  '
    'D6' = '0.9979'
    'E6' = '  -0.15%  '
    'D7' = '0.4947'
    'E7' = '  +4.16%  '
    'D8' = '0.2815'
    'E8' = '  +10.22%  '
    'D9' = '0.06471'
    'E9' = '  +5.83%  '
    'D10' = '1.839.99'
    'E10' = '  +6.91%  '
    'D11' = '16.97'
    'E11' = '  +7.32%  '
    'D12' = '0.07117'
    'E12' = '  +3.25%  '
    'D13' = '0.6629'
    'E13' = '  +11.54%  '
    'D14' = '84.71'
    'E14' = '  +11.13%  '
    'D15' = '4.740'
    'E15' = '  +7.91%  '
    'D16' = '29.785.02'
    'E16' = '  +11.52%  '
    'D17' = '0.9955'
    'E17' = '  -0.47%  '
    'D18' = '0.000007383'
    'E18' = '  +5.72%  '
    'D19' = '12.50'
    'E19' = '  +11.12%  '
    'D20' = '0.9968'
    'E20' = '  -0.25%  '
    'D21' = '2.072.14'
    'E21' = '  +6.80%  '
    'D22' = '4.595'
    'E22' = '  +5.32%  '
    'D23' = '5.445'
    'E23' = '  +7.81%  '
    'D24' = '8.893'
    'E24' = '  +6.96%  '
    'D25' = '143.17'
    'E25' = '  +1.69%  '
    'D26' = '132.50'
    'E26' = '  +25.09%  '
    'D27' = '16.53'
    'E27' = '  +9.47%  '
    'D28' = '1.907'
    'E28' = '  +6.66%  '
    'D29' = '1.403'
    'E29' = '  +2.19%  '
    'D30' = '4.186'
    'E30' = '  +6.41%  '
    'D31' = '0.08559'
    'E31' = '  +8.61%  '
    'D32' = '3.823'
    'E32' = '  +5.44%  '
    'D33' = '0.05004'
    'E33' = '  +8.56%  '
    'D34' = '1.112'
    'E34' = '  +12.09%  '
    'D35' = '0.6803'
    'E35' = '  +11.85%  '
    'D36' = '2.696'
    'E36' = '  +3.90%  '
    'D37' = '2.306'
    'E37' = '  +17.50%  '
    'D38' = '2.733'
    'E38' = '  +9.24%  '
    'D39' = '0.9559'
    'E39' = '  +4.38%  '
    'D40' = '6.151'
    'E40' = '  +8.95%  '
    'D41' = '0.01600'
    'E41' = '  +8.40%  '
    'B42' = 'PaxDollar'
    'C42' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D42' = '0.9966'
    'E42' = '  -0.23%  '
    'B43' = 'Quant'
    'C43' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D43' = '103.43'
    'E43' = '  +3.90%  '
    'D44' = '0.4099'
    'E44' = '  +8.58%  '
    'D45' = '7.268'
    'E45' = '  +8.52%  '
    'D46' = '0.1231'
    'E46' = '  +7.88%  '
    'E47' = '  +4.52%  '
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D48' = '8.196'
    'E48' = '  +6.52%  '
    'B49' = 'Elrond'
    'C49' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D49' = '31.94'
    'E49' = '  +7.66%  '
    'D50' = '1.321'
    'E50' = '  +7.29%  '
    'D51' = '0.3646'
    'E51' = '  +9.78%  '
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
